$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1884057971014493
$ws.Range("C2").Value = 0.5449275362318841
$ws.Range("J2").Value = 0.02318840579710145
$ws.Range("P2").Value = 0.1710144927536232
$ws.Range("S2").Value = 0.07246376811594203
$ws.Range("B3").Value = 0.015625
$ws.Range("C3").Value = 0.02083333333333333
$ws.Range("J3").Value = 0.015625
$ws.Range("P3").Value = 0.7708333333333334
$ws.Range("S3").Value = 0.1770833333333333
$ws.Range("J4").Value = 0.1136363636363636
$ws.Range("P4").Value = 0.5681818181818182
$ws.Range("S4").Value = 0.3181818181818182
$ws.Range("B6").Value = 0.07042253521126761
$ws.Range("D6").Value = 0.004694835680751174
$ws.Range("E6").Value = 0.004694835680751174
$ws.Range("F6").Value = 0.04694835680751173
$ws.Range("J6").Value = 0.3098591549295774
$ws.Range("O6").Value = 0.02347417840375587
$ws.Range("Q6").Value = 0.1549295774647887
$ws.Range("R6").Value = 0.07511737089201878
$ws.Range("S6").Value = 0.3098591549295774
$ws.Range("B7").Value = 0.1288888888888889
$ws.Range("D7").Value = 0.01333333333333333
$ws.Range("F7").Value = 0.04888888888888889
$ws.Range("J7").Value = 0.1022222222222222
$ws.Range("O7").Value = 0.04444444444444445
$ws.Range("Q7").Value = 0.1377777777777778
$ws.Range("R7").Value = 0.08
$ws.Range("S7").Value = 0.4444444444444444
$ws.Range("B8").Value = 0.1085106382978723
$ws.Range("D8").Value = 0.01063829787234043
$ws.Range("F8").Value = 0.06808510638297872
$ws.Range("J8").Value = 0.1042553191489362
$ws.Range("O8").Value = 0.0276595744680851
$ws.Range("Q8").Value = 0.1382978723404255
$ws.Range("R8").Value = 0.1127659574468085
$ws.Range("S8").Value = 0.4297872340425532
$ws.Range("B9").Value = 0.1149425287356322
$ws.Range("D9").Value = 0.01724137931034483
$ws.Range("F9").Value = 0.08620689655172414
$ws.Range("J9").Value = 0.1091954022988506
$ws.Range("O9").Value = 0.02873563218390805
$ws.Range("Q9").Value = 0.1379310344827586
$ws.Range("R9").Value = 0.05172413793103448
$ws.Range("S9").Value = 0.4540229885057471
$ws.Range("B10").Value = 0.1238720262510254
$ws.Range("D10").Value = 0.02871205906480722
$ws.Range("E10").Value = 0.0008203445447087777
$ws.Range("F10").Value = 0.05578342904019688
$ws.Range("J10").Value = 0.1214109926168991
$ws.Range("O10").Value = 0.01886792452830189
$ws.Range("Q10").Value = 0.2042657916324856
$ws.Range("R10").Value = 0.08367514356029532
$ws.Range("S10").Value = 0.3625922887612797
$ws.Range("G11").Value = 0.1621621621621622
$ws.Range("J11").Value = 0.08648648648648649
$ws.Range("K11").Value = 0.2189189189189189
$ws.Range("L11").Value = 0.5135135135135135
$ws.Range("S11").Value = 0.01891891891891892
$ws.Range("G12").Value = 0.6989795918367347
$ws.Range("J12").Value = 0.25
$ws.Range("K12").Value = 0.01020408163265306
$ws.Range("L12").Value = 0.01530612244897959
$ws.Range("S12").Value = 0.02551020408163265
$ws.Range("G13").Value = 0.7333333333333333
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("H15").Value = 0.1084905660377359
$ws.Range("I15").Value = 0.09905660377358491
$ws.Range("J15").Value = 0.3867924528301887
$ws.Range("K15").Value = 0.05188679245283019
$ws.Range("M15").Value = 0.009433962264150943
$ws.Range("O15").Value = 0.03773584905660377
$ws.Range("S15").Value = 0.3066037735849056
$ws.Range("F16").Value = 0.02212389380530973
$ws.Range("H16").Value = 0.1858407079646018
$ws.Range("I16").Value = 0.08849557522123894
$ws.Range("J16").Value = 0.3982300884955752
$ws.Range("K16").Value = 0.1327433628318584
$ws.Range("M16").Value = 0.02212389380530973
$ws.Range("O16").Value = 0.05309734513274336
$ws.Range("S16").Value = 0.09734513274336283
$ws.Range("F17").Value = 0.02551020408163265
$ws.Range("H17").Value = 0.1862244897959184
$ws.Range("I17").Value = 0.09948979591836735
$ws.Range("J17").Value = 0.3826530612244898
$ws.Range("K17").Value = 0.09948979591836735
$ws.Range("M17").Value = 0.01530612244897959
$ws.Range("N17").Value = 0.002551020408163265
$ws.Range("O17").Value = 0.0586734693877551
$ws.Range("S17").Value = 0.1301020408163265
$ws.Range("F18").Value = 0.04020100502512563
$ws.Range("H18").Value = 0.2160804020100502
$ws.Range("I18").Value = 0.09045226130653267
$ws.Range("J18").Value = 0.3768844221105528
$ws.Range("K18").Value = 0.07035175879396985
$ws.Range("M18").Value = 0.02512562814070352
$ws.Range("O18").Value = 0.05527638190954774
$ws.Range("S18").Value = 0.1256281407035176
$ws.Range("F19").Value = 0.02488335925349922
$ws.Range("H19").Value = 0.2262830482115085
$ws.Range("I19").Value = 0.05987558320373251
$ws.Range("J19").Value = 0.3367029548989113
$ws.Range("K19").Value = 0.1438569206842924
$ws.Range("M19").Value = 0.02099533437013997
$ws.Range("O19").Value = 0.05987558320373251
$ws.Range("S19").Value = 0.1275272161741835
